$wb = $excel.ActiveWorkbook

# Updated market-price driven figures (currentAveragePrice*, Leve cost/profit
# columns H:N) per sheet, refreshed by the scheduled market-data runner.

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1374.3334
$ws.Range("I40").Value = 1397.25
$ws.Range("J40").Value = 1356
$ws.Range("K40").Value = 1397.25
$ws.Range("L40").Value = 1356
$ws.Range("M40").Value = -1222.25
$ws.Range("N40").Value = -1706
$ws.Range("H92").Value = 1309.174
$ws.Range("I92").Value = 1155.05
$ws.Range("J92").Value = 2336.6667
$ws.Range("K92").Value = 1155.05
$ws.Range("L92").Value = 2336.6667
$ws.Range("M92").Value = 92.95000000000005
$ws.Range("N92").Value = -4832.6667
$ws.Range("H106").Value = 2572.5
$ws.Range("I106").Value = 1847.3077
$ws.Range("K106").Value = 1847.3077
$ws.Range("M106").Value = -1216.3077
$ws.Range("H112").Value = 1636.8718
$ws.Range("I112").Value = 979.8
$ws.Range("J112").Value = 1733.5
$ws.Range("K112").Value = 2939.4
$ws.Range("L112").Value = 5200.5
$ws.Range("M112").Value = -1831.4
$ws.Range("N112").Value = -7416.5
$ws.Range("H116").Value = 484204.28
$ws.Range("I116").Value = 2003557.6
$ws.Range("J116").Value = 9406.375
$ws.Range("K116").Value = 2003557.6
$ws.Range("L116").Value = 9406.375
$ws.Range("M116").Value = -2000115.6
$ws.Range("N116").Value = -16290.375
$ws.Range("H132").Value = 22668508
$ws.Range("I132").Value = 32581880
$ws.Range("J132").Value = 717468.6
$ws.Range("K132").Value = 97745640
$ws.Range("L132").Value = 2152405.8
$ws.Range("M132").Value = -97743110
$ws.Range("N132").Value = -2157465.8
$ws.Range("H137").Value = 2850.842
$ws.Range("I137").Value = 1341.56
$ws.Range("K137").Value = 4024.68
$ws.Range("M137").Value = -1474.68

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1147.5938
$ws.Range("I61").Value = 890.4483
$ws.Range("J61").Value = 3633.3333
$ws.Range("K61").Value = 890.4483
$ws.Range("L61").Value = 3633.3333
$ws.Range("M61").Value = -678.4483
$ws.Range("N61").Value = -4057.3333
$ws.Range("H74").Value = 3417.5144
$ws.Range("I74").Value = 3310.276
$ws.Range("K74").Value = 3310.276
$ws.Range("M74").Value = -2436.276
$ws.Range("H75").Value = 35000
$ws.Range("J75").Value = 35000
$ws.Range("L75").Value = 35000
$ws.Range("N75").Value = -36748
$ws.Range("H77").Value = 3417.5144
$ws.Range("I77").Value = 3310.276
$ws.Range("K77").Value = 16551.38
$ws.Range("M77").Value = -12183.38
$ws.Range("H78").Value = 35000
$ws.Range("J78").Value = 35000
$ws.Range("L78").Value = 105000
$ws.Range("N78").Value = -113736
$ws.Range("H97").Value = 733.75
$ws.Range("I97").Value = 716.4815
$ws.Range("K97").Value = 716.4815
$ws.Range("M97").Value = -220.4815
$ws.Range("H102").Value = 1528.6
$ws.Range("I102").Value = 1352.9
$ws.Range("J102").Value = 1880
$ws.Range("K102").Value = 1352.9
$ws.Range("L102").Value = 1880
$ws.Range("M102").Value = 269.0999999999999
$ws.Range("N102").Value = -5124
$ws.Range("H136").Value = 1147.5938
$ws.Range("I136").Value = 890.4483
$ws.Range("J136").Value = 3633.3333
$ws.Range("K136").Value = 2671.3449
$ws.Range("L136").Value = 10899.9999
$ws.Range("M136").Value = -121.3449000000001
$ws.Range("N136").Value = -15999.9999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 779.06665
$ws.Range("I94").Value = 811.3077
$ws.Range("J94").Value = 569.5
$ws.Range("K94").Value = 811.3077
$ws.Range("L94").Value = 569.5
$ws.Range("M94").Value = -360.3077
$ws.Range("N94").Value = -1471.5
$ws.Range("H99").Value = 2224.9167
$ws.Range("I99").Value = 1368.1428
$ws.Range("J99").Value = 3424.4
$ws.Range("K99").Value = 1368.1428
$ws.Range("L99").Value = 3424.4
$ws.Range("M99").Value = 129.8571999999999
$ws.Range("N99").Value = -6420.4
$ws.Range("H134").Value = 1403.8909
$ws.Range("I134").Value = 1005.1064
$ws.Range("J134").Value = 3746.75
$ws.Range("K134").Value = 3015.3192
$ws.Range("L134").Value = 11240.25
$ws.Range("M134").Value = -480.3191999999999
$ws.Range("N134").Value = -16310.25

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12502551
$ws.Range("I31").Value = 1411.3478
$ws.Range("J31").Value = 29415858
$ws.Range("K31").Value = 1411.3478
$ws.Range("L31").Value = 29415858
$ws.Range("M31").Value = -1116.3478
$ws.Range("N31").Value = -29416448
$ws.Range("H34").Value = 12502551
$ws.Range("I34").Value = 1411.3478
$ws.Range("J34").Value = 29415858
$ws.Range("K34").Value = 1411.3478
$ws.Range("L34").Value = 29415858
$ws.Range("M34").Value = -1209.3478
$ws.Range("N34").Value = -29416262
$ws.Range("H58").Value = 1667.4819
$ws.Range("I58").Value = 1406.5734
$ws.Range("J58").Value = 4113.5
$ws.Range("K58").Value = 1406.5734
$ws.Range("L58").Value = 4113.5
$ws.Range("M58").Value = -1203.5734
$ws.Range("N58").Value = -4519.5
$ws.Range("H122").Value = 3380.9092
$ws.Range("I122").Value = 1563.3334
$ws.Range("J122").Value = 4062.5
$ws.Range("K122").Value = 4690.0002
$ws.Range("L122").Value = 12187.5
$ws.Range("M122").Value = -2240.0002
$ws.Range("N122").Value = -17087.5
$ws.Range("H132").Value = 1554.2222
$ws.Range("I132").Value = 872.42426
$ws.Range("J132").Value = 3429.1667
$ws.Range("K132").Value = 2617.27278
$ws.Range("L132").Value = 10287.5001
$ws.Range("M132").Value = -87.27278000000024
$ws.Range("N132").Value = -15347.5001
$ws.Range("H134").Value = 2077
$ws.Range("I134").Value = 869.82355
$ws.Range("J134").Value = 3542.8572
$ws.Range("K134").Value = 2609.47065
$ws.Range("L134").Value = 10628.5716
$ws.Range("M134").Value = -74.47064999999975
$ws.Range("N134").Value = -15698.5716
$ws.Range("H136").Value = 1667.4819
$ws.Range("I136").Value = 1406.5734
$ws.Range("J136").Value = 4113.5
$ws.Range("K136").Value = 4219.7202
$ws.Range("L136").Value = 12340.5
$ws.Range("M136").Value = -1669.7202
$ws.Range("N136").Value = -17440.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 655.3333
$ws.Range("I113").Value = 661.5
$ws.Range("K113").Value = 1984.5
$ws.Range("M113").Value = 185.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 735
$ws.Range("I97").Value = 684
$ws.Range("J97").Value = 990
$ws.Range("K97").Value = 684
$ws.Range("L97").Value = 990
$ws.Range("M97").Value = -188
$ws.Range("N97").Value = -1982
$ws.Range("H122").Value = 4024.5557
$ws.Range("I122").Value = 2457.077
$ws.Range("J122").Value = 8100
$ws.Range("K122").Value = 7371.231000000001
$ws.Range("L122").Value = 24300
$ws.Range("M122").Value = -4921.231000000001
$ws.Range("N122").Value = -29200

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4086.1428
$ws.Range("I7").Value = 3437.5
$ws.Range("J7").Value = 4485.3076
$ws.Range("K7").Value = 3437.5
$ws.Range("L7").Value = 4485.3076
$ws.Range("M7").Value = -3325.5
$ws.Range("N7").Value = -4709.3076
$ws.Range("H46").Value = 2375
$ws.Range("I46").Value = 866.6667
$ws.Range("J46").Value = 3280
$ws.Range("K46").Value = 866.6667
$ws.Range("L46").Value = 3280
$ws.Range("M46").Value = -678.6667
$ws.Range("N46").Value = -3656
$ws.Range("H126").Value = 4086.1428
$ws.Range("I126").Value = 3437.5
$ws.Range("J126").Value = 4485.3076
$ws.Range("K126").Value = 10312.5
$ws.Range("L126").Value = 13455.9228
$ws.Range("M126").Value = -7842.5
$ws.Range("N126").Value = -18395.9228
$ws.Range("H132").Value = 9552.578
$ws.Range("I132").Value = 10972.92
$ws.Range("K132").Value = 32918.76
$ws.Range("M132").Value = -30388.76
$ws.Range("H136").Value = 2312.0625
$ws.Range("I136").Value = 1122.5385
$ws.Range("J136").Value = 7466.6665
$ws.Range("K136").Value = 3367.6155
$ws.Range("L136").Value = 22399.9995
$ws.Range("M136").Value = -817.6155000000003
$ws.Range("N136").Value = -27499.9995

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 147816370
$ws.Range("I96").Value = 202121010
$ws.Range("K96").Value = 202121010
$ws.Range("M96").Value = -202119637
$ws.Range("H107").Value = 1101.4
$ws.Range("I107").Value = 1001.3333
$ws.Range("K107").Value = 3003.9999
$ws.Range("M107").Value = -1083.9999
$ws.Range("H132").Value = 4445826.5
$ws.Range("I132").Value = 912.5082
$ws.Range("J132").Value = 23812952
$ws.Range("K132").Value = 2737.5246
$ws.Range("L132").Value = 71438856
$ws.Range("M132").Value = -207.5245999999997
$ws.Range("N132").Value = -71443916
$ws.Range("H136").Value = 1981.7963
$ws.Range("I136").Value = 486.68292
$ws.Range("J136").Value = 6697.154
$ws.Range("K136").Value = 1460.04876
$ws.Range("L136").Value = 20091.462
$ws.Range("M136").Value = 1089.95124
$ws.Range("N136").Value = -25191.462
